$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume data per latest scrape
$ws.Range('D2').Value = '29.885.18'
$ws.Range('E2').Value = '  +0.78%  '
$ws.Range('D3').Value = '1.628.20'
$ws.Range('E3').Value = '  +1.07%  '
$ws.Range('E4').Value = '  +0.65%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.30'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.82%  '
$ws.Range('E6').Value = '  -0.57%  '
$ws.Range('E7').Value = '  +0.66%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '28.50'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.02%  '
$ws.Range('E9').Value = '  -0.23%  '
$ws.Range('E10').Value = '  -0.09%  '
$ws.Range('D11').Value = '0.0901'
$ws.Range('E11').Value = '  -0.75%  '
$ws.Range('D12').Value = '1.856.10'
$ws.Range('E12').Value = '  +0.81%  '
$ws.Range('D13').Value = '1.627.83'
$ws.Range('E13').Value = '  +1.15%  '
$ws.Range('E14').Value = '  -1.49%  '
$ws.Range('D15').Value = '9.23'
$ws.Range('E15').Value = '  +6.55%  '
$ws.Range('D16').Value = '29.912.01'
$ws.Range('E16').Value = '  +0.83%  '
$ws.Range('E17').Value = '  -0.46%  '
$ws.Range('D18').Value = '64.02'
$ws.Range('E18').Value = '  -0.98%  '
$ws.Range('D19').Value = '240.53'
$ws.Range('E19').Value = '  +0.04%  '
$ws.Range('D20').Value = '0.0₃0701'
$ws.Range('E20').Value = '  -0.47%  '
$ws.Range('E21').Value = '  +0.50%  '
$ws.Range('E22').Value = '  +1.06%  '
$ws.Range('E23').Value = '  +1.58%  '
$ws.Range('D24').Value = '2.16'
$ws.Range('E24').Value = '  +1.60%  '
$ws.Range('D25').Value = '157.81'
$ws.Range('E25').Value = '  +0.78%  '
$ws.Range('D26').Value = '15.44'
$ws.Range('E26').Value = '  -0.93%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.110'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -0.46%  '
$ws.Range('D28').Value = '6.53'
$ws.Range('E28').Value = '  -0.39%  '
$ws.Range('E29').Value = '  +0.59%  '
$ws.Range('E30').Value = '  +1.17%  '
$ws.Range('E31').Value = '  +3.04%  '
$ws.Range('E32').Value = '  +1.77%  '
$ws.Range('D33').Value = '3.17'
$ws.Range('D34').Value = '1.425.20'
$ws.Range('E34').Value = '  -1.10%  '
$ws.Range('E35').Value = '  +4.29%  '
$ws.Range('D36').Value = '1.02'
$ws.Range('E36').Value = '  -2.17%  '
$ws.Range('E37').Value = '  -4.67%  '
$ws.Range('E38').Value = '  +0.24%  '
$ws.Range('E39').Value = '  -0.03%  '
$ws.Range('B40').Value = 'ImmutableX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D40').Value = '0.555'
$ws.Range('E40').Value = '  +0.18%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').Value = '74.72'
$ws.Range('E41').Value = '  +7.08%  '
$ws.Range('B42').Value = 'Kaspa'
$ws.Range('C42').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D42').Value = '0.0498'
$ws.Range('E42').Value = '  -1.51%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = '1.98'
$ws.Range('E43').Value = '  -1.81%  '
$ws.Range('D44').Value = '0.825'
$ws.Range('E44').Value = '  +0.19%  '
$ws.Range('E45').Value = '  +0.65%  '
$ws.Range('E46').Value = '  +0.32%  '
$ws.Range('D47').Value = '5.33'
$ws.Range('E47').Value = '  -2.27%  '
$ws.Range('D48').Value = '1.768.07'
$ws.Range('E48').Value = '  +1.01%  '
$ws.Range('D49').Value = '49.11'
$ws.Range('E49').Value = '  -9.21%  '
$ws.Range('D50').Value = '90.79'
$ws.Range('E50').Value = '  +3.83%  '
$ws.Range('E51').Value = '  +8.90%  '
